$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Select the whole row about to be removed (mirrors the original authoring
# session: user selects row 6 - the "5. Selectionner un village de 1ere
# ligne" question - then deletes it, shifting all following rows up).
$survey.Activate()
$survey.Range("A6:XFD6").Select()
$survey.Rows.Item(6).Delete()

# Update the form title / form id on the settings sheet for the new
# version (V3 -> V3.1).
$settings.Cells.Item(2, 1).Value = "(Sept 2023) ONCHO Pre Stop - 2. Participants V3.1"
$settings.Cells.Item(2, 2).Value = "civ_oncho_ia_202309_2_participant_v3_1"

# The settings tab becomes the active / selected tab, with B2 selected.
$settings.Activate()
$settings.Range("B2").Select()
